# Generate Report for Handback
# Renames the two handed-back file UUIDs (and the regenerated .xlf hash) that
# are referenced throughout the Overview / zh-cn / de-de sheets, and updates
# the corresponding handoff/handback timestamps.

$wb = $excel.ActiveWorkbook

$oldUuid1 = "362a6b95-1b89-4d8f-95ee-f7f51efca2d1"
$newUuid1 = "aac1a086-635d-4bc7-8d61-23c24bef0bb6"
$oldUuid2 = "4d150459-f398-47c5-8562-1c426fe27a5a"
$newUuid2 = "ffff0dc19648-959a-4389-8364-2c5b86ef7bed"

$newMd1 = "$newUuid1.md"
$newMd2 = "$newUuid2.md"

$newHash = "c31798f111cbfeadbf8bd1000277a568068217a3"
$newXlfZh = "$newUuid1.$newHash.zh-cn.xlf"
$newXlfDe = "$newUuid1.$newHash.de-de.xlf"

$newHandoffZh = "2016-03-24 05:10:46"
$newHandbackZh = "2016-03-24 05:11:09"
$newHandoffDe = "2016-03-24 05:10:50"
$newHandbackDe = "2016-03-24 05:11:16"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd1
$wsOverview.Range("A3").Value = $newMd2

foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $newMd1
    }
    if ($addr -eq '$A$3') {
        $h.TextToDisplay = $newMd2
    }
}

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMd1
$wsZh.Range("D2").Value = $newXlfZh
$wsZh.Range("E2").Value = $newHandoffZh
$wsZh.Range("F2").Value = $newMd1
$wsZh.Range("G2").Value = $newXlfZh
$wsZh.Range("H2").Value = $newHandbackZh

$wsZh.Range("A3").Value = $newMd2
$wsZh.Range("D3").Value = $newXlfZh
$wsZh.Range("E3").Value = $newHandoffZh
$wsZh.Range("F3").Value = $newMd2
$wsZh.Range("G3").Value = $newXlfZh
$wsZh.Range("H3").Value = $newHandbackZh

foreach ($h in $wsZh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') { $h.TextToDisplay = $newMd1 }
    if ($addr -eq '$D$2') { $h.TextToDisplay = $newXlfZh }
    if ($addr -eq '$F$2') { $h.TextToDisplay = $newMd1 }
    if ($addr -eq '$G$2') { $h.TextToDisplay = $newXlfZh }
    if ($addr -eq '$A$3') { $h.TextToDisplay = $newMd2 }
    if ($addr -eq '$D$3') { $h.TextToDisplay = $newXlfZh }
    if ($addr -eq '$F$3') { $h.TextToDisplay = $newMd2 }
    if ($addr -eq '$G$3') { $h.TextToDisplay = $newXlfZh }
}

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMd1
$wsDe.Range("D2").Value = $newXlfDe
$wsDe.Range("E2").Value = $newHandoffDe
$wsDe.Range("F2").Value = $newMd1
$wsDe.Range("G2").Value = $newXlfDe
$wsDe.Range("H2").Value = $newHandbackDe

$wsDe.Range("A3").Value = $newMd2
$wsDe.Range("D3").Value = $newXlfDe
$wsDe.Range("E3").Value = $newHandoffDe
$wsDe.Range("F3").Value = $newMd2
$wsDe.Range("G3").Value = $newXlfDe
$wsDe.Range("H3").Value = $newHandbackDe

foreach ($h in $wsDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') { $h.TextToDisplay = $newMd1 }
    if ($addr -eq '$D$2') { $h.TextToDisplay = $newXlfDe }
    if ($addr -eq '$F$2') { $h.TextToDisplay = $newMd1 }
    if ($addr -eq '$G$2') { $h.TextToDisplay = $newXlfDe }
    if ($addr -eq '$A$3') { $h.TextToDisplay = $newMd2 }
    if ($addr -eq '$D$3') { $h.TextToDisplay = $newXlfDe }
    if ($addr -eq '$F$3') { $h.TextToDisplay = $newMd2 }
    if ($addr -eq '$G$3') { $h.TextToDisplay = $newXlfDe }
}
